{"js": "// Replace the arithmetic-problem text in every cell of the first\n// table in the document body with the new values from the commit,\n// preserving row/column order (row-major, 20 rows x 5 columns).\nconst newValues = [\n  [\"54+37=\", \"87+5=\", \"70-46=\", \"59+23=\", \"15+38=\"],\n  [\"94-66=\", \"4+67=\", \"51-34=\", \"35+26=\", \"46+17=\"],\n  [\"34-5=\", \"25-16=\", \"47+29=\", \"76+16=\", \"17+34=\"],\n  [\"48+45=\", \"53-7=\", \"39+7=\", \"27+47=\", \"7+65=\"],\n  [\"87+4=\", \"84-55=\", \"83-29=\", \"64+9=\", \"47+28=\"],\n  [\"31-17=\", \"72-58=\", \"93-65=\", \"63-59=\", \"92-8=\"],\n  [\"61-49=\", \"66-9=\", \"91-84=\", \"29+7=\", \"89+9=\"],\n  [\"70-46=\", \"67+9=\", \"47+8=\", \"38+48=\", \"62-17=\"],\n  [\"66-37=\", \"61-42=\", \"26+67=\", \"92-87=\", \"80-41=\"],\n  [\"92-77=\", \"49+9=\", \"92-45=\", \"24+27=\", \"21-6=\"],\n  [\"94-27=\", \"42-27=\", \"53+39=\", \"59+2=\", \"9+88=\"],\n  [\"15+66=\", \"97-8=\", \"35+58=\", \"53-34=\", \"83-39=\"],\n  [\"59+34=\", \"28+14=\", \"43-27=\", \"54-16=\", \"93-17=\"],\n  [\"80-36=\", \"62-56=\", \"47+17=\", \"70-13=\", \"54-19=\"],\n  [\"48+23=\", \"14+68=\", \"81-44=\", \"49+33=\", \"81-38=\"],\n  [\"33+28=\", \"7+75=\", \"60-56=\", \"49+49=\", \"74+9=\"],\n  [\"59+34=\", \"14+19=\", \"86-38=\", \"23+38=\", \"16+68=\"],\n  [\"6+39=\", \"96-29=\", \"38+25=\", \"82-78=\", \"15+29=\"],\n  [\"15+16=\", \"68+6=\", \"96-28=\", \"4+78=\", \"93-44=\"],\n  [\"86-7=\", \"45+6=\", \"43-26=\", \"48+49=\", \"90-52=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst expectedRows = newValues.length;\nconst expectedCols = newValues[0].length;\nif (table.rowCount !== expectedRows || table.values[0].length !== expectedCols) {\n  throw new Error(\n    `Unexpected table shape: got ${table.rowCount}x${table.values[0].length}, ` +\n    `expected ${expectedRows}x${expectedCols}`\n  );\n}\n\n// Bulk-write every cell's text in one shot, preserving each cell's\n// existing run/paragraph formatting (only the text content changes).\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# Replace the arithmetic-problem text in every cell of the first\n# table in the document with the new values from the commit,\n# preserving row/column order (row-major, 20 rows x 5 columns).\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$newValues = @(\n  \"54+37=\",\n  \"87+5=\",\n  \"70-46=\",\n  \"59+23=\",\n  \"15+38=\",\n  \"94-66=\",\n  \"4+67=\",\n  \"51-34=\",\n  \"35+26=\",\n  \"46+17=\",\n  \"34-5=\",\n  \"25-16=\",\n  \"47+29=\",\n  \"76+16=\",\n  \"17+34=\",\n  \"48+45=\",\n  \"53-7=\",\n  \"39+7=\",\n  \"27+47=\",\n  \"7+65=\",\n  \"87+4=\",\n  \"84-55=\",\n  \"83-29=\",\n  \"64+9=\",\n  \"47+28=\",\n  \"31-17=\",\n  \"72-58=\",\n  \"93-65=\",\n  \"63-59=\",\n  \"92-8=\",\n  \"61-49=\",\n  \"66-9=\",\n  \"91-84=\",\n  \"29+7=\",\n  \"89+9=\",\n  \"70-46=\",\n  \"67+9=\",\n  \"47+8=\",\n  \"38+48=\",\n  \"62-17=\",\n  \"66-37=\",\n  \"61-42=\",\n  \"26+67=\",\n  \"92-87=\",\n  \"80-41=\",\n  \"92-77=\",\n  \"49+9=\",\n  \"92-45=\",\n  \"24+27=\",\n  \"21-6=\",\n  \"94-27=\",\n  \"42-27=\",\n  \"53+39=\",\n  \"59+2=\",\n  \"9+88=\",\n  \"15+66=\",\n  \"97-8=\",\n  \"35+58=\",\n  \"53-34=\",\n  \"83-39=\",\n  \"59+34=\",\n  \"28+14=\",\n  \"43-27=\",\n  \"54-16=\",\n  \"93-17=\",\n  \"80-36=\",\n  \"62-56=\",\n  \"47+17=\",\n  \"70-13=\",\n  \"54-19=\",\n  \"48+23=\",\n  \"14+68=\",\n  \"81-44=\",\n  \"49+33=\",\n  \"81-38=\",\n  \"33+28=\",\n  \"7+75=\",\n  \"60-56=\",\n  \"49+49=\",\n  \"74+9=\",\n  \"59+34=\",\n  \"14+19=\",\n  \"86-38=\",\n  \"23+38=\",\n  \"16+68=\",\n  \"6+39=\",\n  \"96-29=\",\n  \"38+25=\",\n  \"82-78=\",\n  \"15+29=\",\n  \"15+16=\",\n  \"68+6=\",\n  \"96-28=\",\n  \"4+78=\",\n  \"93-44=\",\n  \"86-7=\",\n  \"45+6=\",\n  \"43-26=\",\n  \"48+49=\",\n  \"90-52=\"\n)\n\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\nif (($rows * $cols) -ne $newValues.Count) {\n  throw \"Unexpected table shape: $rows x $cols cells, expected $($newValues.Count) values\"\n}\n\n# Cell.Range.Text replaces only the cell's text run, leaving the\n# existing paragraph/run formatting (font, size, alignment) intact.\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $tbl.Cell($r, $c).Range.Text = $newValues[$idx]\n    $idx = $idx + 1\n  }\n}\n"}
